# Apply weekly crime-data refresh to the CompStat 33rd Precinct report.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Header text updates (Volume/Number, and the "Report Covering the Week" line)
# ---------------------------------------------------------------------------
$ws.Range("A8").Replace("26", "27")
$ws.Range("C9").Replace("6/24/2024", "7/1/2024")
$ws.Range("C9").Replace("6/30/2024", "7/7/2024")

# ---------------------------------------------------------------------------
# Helper: turn a numeric cell into a "N/A"-style text cell ("0" or "***.*"),
# while keeping the exact same look (style) as a neighboring text cell that
# already uses that formatting.
# ---------------------------------------------------------------------------
function Set-TextCell($cellRef, $text, $formatSourceRef) {
    $ws.Range($cellRef).Value = "'" + $text
    $ws.Range($formatSourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
}

# Helper: turn a text "N/A"-style cell into a normal numeric cell, keeping the
# same look (style) as a neighboring numeric cell.
function Set-NumericCell($cellRef, $value, $formatSourceRef) {
    $ws.Range($formatSourceRef).Copy()
    $ws.Range($cellRef).PasteSpecial(-4122)
    $ws.Range($cellRef).Value = $value
}

# ---------------------------------------------------------------------------
# Row 14 - Murder
# ---------------------------------------------------------------------------
Set-TextCell "D14" "0" "C14"
Set-TextCell "E14" "***.*" "N14"
$ws.Range("F14").Value = 1
$ws.Range("H14").Value = 0

# ---------------------------------------------------------------------------
# Row 16 - Robbery
# ---------------------------------------------------------------------------
$ws.Range("D16").Value = 3
$ws.Range("E16").Value = 0
$ws.Range("F16").Value = 8
$ws.Range("G16").Value = 18
$ws.Range("H16").Value = -55.555555555555
$ws.Range("I16").Value = 92
$ws.Range("J16").Value = 75
$ws.Range("K16").Value = 22.666666666666
$ws.Range("L16").Value = -4.166666666666
$ws.Range("M16").Value = -24.590163934426

# ---------------------------------------------------------------------------
# Row 17 - Fel. Assault
# ---------------------------------------------------------------------------
$ws.Range("C17").Value = 5
$ws.Range("D17").Value = 2
$ws.Range("E17").Value = 150
$ws.Range("F17").Value = 21
$ws.Range("G17").Value = 21
$ws.Range("H17").Value = 0
$ws.Range("I17").Value = 132
$ws.Range("J17").Value = 123
$ws.Range("K17").Value = 7.317073170731
$ws.Range("L17").Value = -7.042253521126
$ws.Range("M17").Value = 51.724137931034

# ---------------------------------------------------------------------------
# Row 18 - Burglary
# ---------------------------------------------------------------------------
$ws.Range("D18").Value = 1
$ws.Range("F18").Value = 3
$ws.Range("G18").Value = 8
$ws.Range("H18").Value = -62.5
$ws.Range("J18").Value = 70
$ws.Range("K18").Value = -27.142857142857
$ws.Range("L18").Value = -5.555555555555
$ws.Range("M18").Value = 2

# ---------------------------------------------------------------------------
# Row 19 - Gr. Larceny
# ---------------------------------------------------------------------------
$ws.Range("D19").Value = 8
$ws.Range("E19").Value = 25
$ws.Range("F19").Value = 34
$ws.Range("G19").Value = 33
$ws.Range("H19").Value = 3.030303030303
$ws.Range("I19").Value = 187
$ws.Range("J19").Value = 167
$ws.Range("K19").Value = 11.976047904191
$ws.Range("L19").Value = 17.610062893081
$ws.Range("M19").Value = 36.496350364963

# ---------------------------------------------------------------------------
# Row 20 - G.L.A.
# ---------------------------------------------------------------------------
$ws.Range("C20").Value = 2
$ws.Range("D20").Value = 2
$ws.Range("E20").Value = 0
$ws.Range("F20").Value = 13
$ws.Range("G20").Value = 15
$ws.Range("H20").Value = -13.333333333333
$ws.Range("I20").Value = 55
$ws.Range("J20").Value = 74
$ws.Range("K20").Value = -25.675675675675
$ws.Range("L20").Value = -31.25
$ws.Range("M20").Value = 57.142857142857

# ---------------------------------------------------------------------------
# Row 21 - TOTAL
# ---------------------------------------------------------------------------
$ws.Range("D21").Value = 16
$ws.Range("E21").Value = 25
$ws.Range("F21").Value = 81
$ws.Range("G21").Value = 96
$ws.Range("H21").Value = -15.625
$ws.Range("I21").Value = 526
$ws.Range("J21").Value = 517
$ws.Range("K21").Value = 1.740812379110
$ws.Range("L21").Value = -2.411873840445
$ws.Range("M21").Value = 18.468468468468

# ---------------------------------------------------------------------------
# Row 22 - Transit
# ---------------------------------------------------------------------------
Set-NumericCell "C22" 1 "F22"
$ws.Range("F22").Value = 2
$ws.Range("H22").Value = 100
$ws.Range("I22").Value = 7
$ws.Range("K22").Value = 0
$ws.Range("L22").Value = -56.25
$ws.Range("M22").Value = 16.666666666666

# ---------------------------------------------------------------------------
# Row 23 - Housing
# ---------------------------------------------------------------------------
Set-TextCell "C23" "0" "D23"
$ws.Range("F23").Value = 2
$ws.Range("H23").Value = 100

# ---------------------------------------------------------------------------
# Row 24 - Petit Larceny
# ---------------------------------------------------------------------------
$ws.Range("C24").Value = 17
$ws.Range("D24").Value = 18
$ws.Range("E24").Value = -5.555555555555
$ws.Range("F24").Value = 68
$ws.Range("G24").Value = 73
$ws.Range("H24").Value = -6.849315068493
$ws.Range("I24").Value = 461
$ws.Range("J24").Value = 488
$ws.Range("K24").Value = -5.532786885245
$ws.Range("L24").Value = -41.348600508905
$ws.Range("M24").Value = 62.323943661971

# ---------------------------------------------------------------------------
# Row 25 - Retail Theft
# ---------------------------------------------------------------------------
$ws.Range("D25").Value = 4
$ws.Range("E25").Value = -25
$ws.Range("F25").Value = 31
$ws.Range("G25").Value = 21
$ws.Range("H25").Value = 47.619047619047
$ws.Range("I25").Value = 137
$ws.Range("J25").Value = 127
$ws.Range("K25").Value = 7.874015748031
$ws.Range("L25").Value = -71.868583162217

# ---------------------------------------------------------------------------
# Row 26 - Misd. Assault
# ---------------------------------------------------------------------------
$ws.Range("C26").Value = 6
$ws.Range("E26").Value = -45.454545454545
$ws.Range("F26").Value = 34
$ws.Range("G26").Value = 38
$ws.Range("H26").Value = -10.526315789473
$ws.Range("I26").Value = 174
$ws.Range("J26").Value = 213
$ws.Range("K26").Value = -18.309859154929
$ws.Range("L26").Value = -20.183486238532
$ws.Range("M26").Value = -30.4

# ---------------------------------------------------------------------------
# Row 28 - Other Sex Crimes
# ---------------------------------------------------------------------------
$ws.Range("C28").Value = 1
Set-TextCell "D28" "0" "A28"
Set-TextCell "E28" "***.*" "M28"
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 23
$ws.Range("K28").Value = 27.777777777777
$ws.Range("L28").Value = -17.857142857142
